{"js": "// Update the quantitative-genetics summary table:\n//   \"Date of first flower\"   : h2   0.000 -> 0.039\n//                               CVA  0.001 -> 0.010\n//   \"Date of first follicle\" : h2   0.001 -> 0.235\n//                               QST  0.000 -> 0.005\n//                               CVA  0.000 -> 0.034\n// (QST for \"Date of first flower\" stays 0.000 and is left untouched.)\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Map of row label -> { columnIndex(0-based, 1=h2,2=QST,3=CVA): newValue }\nconst updates = {\n  \"Date of first flower\": { 1: \"0.039\", 3: \"0.010\" },\n  \"Date of first follicle\": { 1: \"0.235\", 2: \"0.005\", 3: \"0.034\" }\n};\n\nconst rows = table.values;\nfor (let r = 0; r < rows.length; r++) {\n  const label = rows[r][0];\n  if (Object.prototype.hasOwnProperty.call(updates, label)) {\n    const cols = updates[label];\n    for (const colIndex of Object.keys(cols)) {\n      const cell = table.getCell(r, Number(colIndex));\n      cell.value = cols[colIndex];\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the quantitative-genetics summary table:\n#   \"Date of first flower\"   : h2   0.000 -> 0.039\n#                               CVA  0.001 -> 0.010\n#   \"Date of first follicle\" : h2   0.001 -> 0.235\n#                               QST  0.000 -> 0.005\n#                               CVA  0.000 -> 0.034\n# (QST for \"Date of first flower\" stays 0.000 and is left untouched.)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $label = $t.Cell($r, 1).Range.Text\n\n    if ($label -like \"Date of first flower*\") {\n        $t.Cell($r, 2).Range.Text = \"0.039\"\n        $t.Cell($r, 4).Range.Text = \"0.010\"\n    }\n    elseif ($label -like \"Date of first follicle*\") {\n        $t.Cell($r, 2).Range.Text = \"0.235\"\n        $t.Cell($r, 3).Range.Text = \"0.005\"\n        $t.Cell($r, 4).Range.Text = \"0.034\"\n    }\n}\n"}
